$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.747.32"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.906.79"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "239.44"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "0.4937"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "0.2963"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "0.06739"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "1.917.66"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "17.05"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").Value = "0.07350"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "5.169"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "88.36"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "0.6707"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "30.707.44"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "0.000007910"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "13.49"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "2.159.80"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "5.342"
$ws.Range("E21").Value = "  +11.23%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "194.03"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").Value = "6.261"
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").Value = "9.617"
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("D26").Value = "163.18"
$ws.Range("E26").Value = "  +3.45%  "
$ws.Range("D27").Value = "18.60"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "1.942"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").Value = "1.475"
$ws.Range("E29").Value = "  +5.28%  "
$ws.Range("D30").Value = "4.453"
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("D31").Value = "0.09148"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "4.041"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "0.05264"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "0.7431"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").Value = "1.109"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").Value = "2.730"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").Value = "0.01826"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "2.711"
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("D39").Value = "0.9217"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "2.079"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").Value = "74.89"
$ws.Range("E41").Value = "  +29.57%  "
$ws.Range("D42").Value = "0.4441"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "5.951"
$ws.Range("E43").Value = "  +3.92%  "
$ws.Range("D44").Value = "106.80"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").Value = "0.9991"
$ws.Range("D46").Value = "0.1388"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("D47").Value = "7.539"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "35.47"
$ws.Range("E48").Value = "  +4.95%  "
$ws.Range("D49").Value = "9.069"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("D50").Value = "0.05856"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "0.3997"
$ws.Range("E51").Value = "  +1.16%  "
